$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 6284
$ws.Range("I100").Value = 1916.6666
$ws.Range("J100").Value = 11524.8
$ws.Range("K100").Value = 1916.6666
$ws.Range("L100").Value = 11524.8
$ws.Range("M100").Value = -1375.6666
$ws.Range("N100").Value = -12606.8
$ws.Range("H111").Value = 3941.111
$ws.Range("I111").Value = 3713.3333
$ws.Range("J111").Value = 4396.6665
$ws.Range("K111").Value = 11139.9999
$ws.Range("L111").Value = 13189.9995
$ws.Range("M111").Value = -8072.999899999999
$ws.Range("N111").Value = -19323.9995
$ws.Range("H132").Value = 487753.1
$ws.Range("I132").Value = 546223.5
$ws.Range("K132").Value = 1638670.5
$ws.Range("M132").Value = -1636140.5
$ws.Range("H135").Value = 1573.037
$ws.Range("I135").Value = 1597.409
$ws.Range("K135").Value = 14376.681
$ws.Range("M135").Value = -11841.681
$ws.Range("H136").Value = 89994.5
$ws.Range("J136").Value = 89994.5
$ws.Range("L136").Value = 89994.5
$ws.Range("N136").Value = -100194.5
$ws.Range("H137").Value = 2168.4583
$ws.Range("I137").Value = 1499.421
$ws.Range("J137").Value = 2606.7932
$ws.Range("K137").Value = 4498.263
$ws.Range("L137").Value = 7820.3796
$ws.Range("M137").Value = -1948.263
$ws.Range("N137").Value = -12920.3796
$ws.Range("H138").Value = 3473.8833
$ws.Range("I138").Value = 2787.4
$ws.Range("J138").Value = 3702.7112
$ws.Range("K138").Value = 8362.200000000001
$ws.Range("L138").Value = 11108.1336
$ws.Range("M138").Value = -3222.200000000001
$ws.Range("N138").Value = -21388.1336
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6858677.5
$ws.Range("I32").Value = 9621176
$ws.Range("K32").Value = 9621176
$ws.Range("M32").Value = -9620889
$ws.Range("H55").Value = 99999.28999999999
$ws.Range("J55").Value = 99999.28999999999
$ws.Range("L55").Value = 99999.28999999999
$ws.Range("N55").Value = -100629.29
$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676
$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340
$ws.Range("H102").Value = 4649.1113
$ws.Range("I102").Value = 4579.1333
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 4579.1333
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = -2957.1333
$ws.Range("N102").Value = -8243
$ws.Range("H105").Value = 70370
$ws.Range("J105").Value = 70370
$ws.Range("L105").Value = 70370
$ws.Range("N105").Value = -77358
$ws.Range("H110").Value = 2931.5789
$ws.Range("I110").Value = 2282.4119
$ws.Range("K110").Value = 2282.4119
$ws.Range("M110").Value = -237.4119000000001
$ws.Range("H112").Value = 53960
$ws.Range("J112").Value = 59900
$ws.Range("L112").Value = 59900
$ws.Range("N112").Value = -62854
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 3666.6667
$ws.Range("K122").Value = 11000.0001
$ws.Range("M122").Value = -8550.000100000001
$ws.Range("H134").Value = 87359.336
$ws.Range("J134").Value = 87359.336
$ws.Range("L134").Value = 87359.336
$ws.Range("N134").Value = -97499.336
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 93498.336
$ws.Range("J35").Value = 93498.336
$ws.Range("L35").Value = 93498.336
$ws.Range("N35").Value = -94118.336
$ws.Range("H105").Value = 3349.5173
$ws.Range("I105").Value = 3262.2354
$ws.Range("J105").Value = 3473.1667
$ws.Range("K105").Value = 3262.2354
$ws.Range("L105").Value = 3473.1667
$ws.Range("M105").Value = -1515.2354
$ws.Range("N105").Value = -6967.1667
$ws.Range("H131").Value = 75000
$ws.Range("J131").Value = 75000
$ws.Range("L131").Value = 75000
$ws.Range("N131").Value = -85080
$ws.Range("H133").Value = 81609
$ws.Range("I133").Value = 81609
$ws.Range("K133").Value = 81609
$ws.Range("M133").Value = -76549
$ws.Range("H139").Value = 64991.25
$ws.Range("J139").Value = 64991.25
$ws.Range("L139").Value = 64991.25
$ws.Range("N139").Value = -75271.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7523.1665
$ws.Range("I31").Value = 1529.0625
$ws.Range("J31").Value = 10047
$ws.Range("K31").Value = 1529.0625
$ws.Range("L31").Value = 10047
$ws.Range("M31").Value = -1234.0625
$ws.Range("N31").Value = -10637
$ws.Range("H34").Value = 7523.1665
$ws.Range("I34").Value = 1529.0625
$ws.Range("J34").Value = 10047
$ws.Range("K34").Value = 1529.0625
$ws.Range("L34").Value = 10047
$ws.Range("M34").Value = -1327.0625
$ws.Range("N34").Value = -10451
$ws.Range("H55").Value = 8399.75
$ws.Range("I55").Value = 8399.75
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 8399.75
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -8084.75
$ws.Range("N55").ClearContents()
$ws.Range("H62").Value = 4480
$ws.Range("I62").Value = 4286.5
$ws.Range("K62").Value = 4286.5
$ws.Range("M62").Value = -3662.5
$ws.Range("H65").Value = 4480
$ws.Range("I65").Value = 4286.5
$ws.Range("K65").Value = 21432.5
$ws.Range("M65").Value = -18312.5
$ws.Range("H107").Value = 951.2353000000001
$ws.Range("I107").Value = 726.04
$ws.Range("K107").Value = 726.04
$ws.Range("M107").Value = 1193.96
$ws.Range("H134").Value = 4383.3657
$ws.Range("I134").Value = 1518.4103
$ws.Range("J134").Value = 60250
$ws.Range("K134").Value = 4555.2309
$ws.Range("L134").Value = 180750
$ws.Range("M134").Value = -2020.2309
$ws.Range("N134").Value = -185820
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185274.44
$ws.Range("I2").Value = 312533.38
$ws.Range("K2").Value = 1875200.28
$ws.Range("M2").Value = -1875087.28
$ws.Range("H68").Value = 8320.866
$ws.Range("I68").Value = 2983
$ws.Range("J68").Value = 9655.333000000001
$ws.Range("K68").Value = 8949
$ws.Range("L68").Value = 28965.999
$ws.Range("M68").Value = -8138
$ws.Range("N68").Value = -30587.999
$ws.Range("H71").Value = 8320.866
$ws.Range("I71").Value = 2983
$ws.Range("J71").Value = 9655.333000000001
$ws.Range("K71").Value = 26847
$ws.Range("L71").Value = 86897.997
$ws.Range("M71").Value = -22791
$ws.Range("N71").Value = -95009.997
$ws.Range("H102").Value = 8067.222
$ws.Range("J102").Value = 9515.143
$ws.Range("L102").Value = 28545.429
$ws.Range("N102").Value = -33413.429
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 2351.7
$ws.Range("I114").Value = 547.6923
$ws.Range("J114").Value = 5702
$ws.Range("K114").Value = 1643.0769
$ws.Range("L114").Value = 17106
$ws.Range("M114").Value = 1610.9231
$ws.Range("N114").Value = -23614
$ws.Range("H131").Value = 8992.933999999999
$ws.Range("I131").Value = 881.55554
$ws.Range("J131").Value = 12469.238
$ws.Range("K131").Value = 2644.66662
$ws.Range("L131").Value = 37407.714
$ws.Range("M131").Value = 2395.33338
$ws.Range("N131").Value = -47487.714
$ws.Range("H132").Value = 3303.9583
$ws.Range("I132").Value = 2730.3
$ws.Range("K132").Value = 24572.7
$ws.Range("M132").Value = -22042.7
$ws.Range("H137").Value = 1821.2084
$ws.Range("I137").Value = 1942.2142
$ws.Range("J137").Value = 1651.8
$ws.Range("K137").Value = 5826.642599999999
$ws.Range("L137").Value = 4955.4
$ws.Range("M137").Value = -726.6425999999992
$ws.Range("N137").Value = -15155.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 45000
$ws.Range("J52").Value = 45000
$ws.Range("L52").Value = 45000
$ws.Range("N52").Value = -45518
$ws.Range("H103").Value = 34650.5
$ws.Range("J103").Value = 34650.5
$ws.Range("L103").Value = 34650.5
$ws.Range("N103").Value = -36994.5
$ws.Range("H111").Value = 59999
$ws.Range("J111").Value = 59999
$ws.Range("L111").Value = 59999
$ws.Range("N111").Value = -66133
$ws.Range("H132").Value = 2494
$ws.Range("I132").Value = 2597.375
$ws.Range("K132").Value = 7792.125
$ws.Range("M132").Value = -5262.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 94181.82000000001
$ws.Range("J22").Value = 3600
$ws.Range("L22").Value = 3600
$ws.Range("N22").Value = -4190
$ws.Range("H27").Value = 94181.82000000001
$ws.Range("J27").Value = 3600
$ws.Range("L27").Value = 3600
$ws.Range("N27").Value = -3814
$ws.Range("H45").Value = 19041
$ws.Range("I45").Value = 19041
$ws.Range("K45").Value = 19041
$ws.Range("M45").Value = -18634
$ws.Range("H46").Value = 1608.091
$ws.Range("I46").Value = 1222.5
$ws.Range("J46").Value = 1828.4286
$ws.Range("K46").Value = 1222.5
$ws.Range("L46").Value = 1828.4286
$ws.Range("M46").Value = -1034.5
$ws.Range("N46").Value = -2204.4286
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H69").Value = 78942
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 78942
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 78942
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -80564
$ws.Range("H72").Value = 78942
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 78942
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 236826
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -244938
$ws.Range("H122").Value = 4383.0835
$ws.Range("I122").Value = 3969.7
$ws.Range("J122").Value = 6450
$ws.Range("K122").Value = 11909.1
$ws.Range("L122").Value = 19350
$ws.Range("M122").Value = -9459.099999999999
$ws.Range("N122").Value = -24250
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3037.5
$ws.Range("I126").Value = 3171
$ws.Range("K126").Value = 9513
$ws.Range("M126").Value = -7043
